$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "123" and "1" look like numbers, so a plain .Value assignment would be
# auto-coerced to a number. Write them as text-producing formulas instead,
# then convert the formula cells in place to plain values (Copy +
# PasteSpecial values-only). That keeps the literal text (and its shared
# string) without ever touching NumberFormat/Style (which would otherwise
# register a stray cell style in styles.xml).
$ws.Range("A2").Value = "31-12-2024"
$ws.Range("B2").Formula = '="123"'
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)
$ws.Range("C2").Formula = '="1"'
$ws.Range("C2").Copy()
$ws.Range("C2").PasteSpecial(-4163)
$ws.Range("D2").Value = "anadora"
$ws.Range("E2").Value = "ana dora"
$ws.Range("F2").Formula = '="1"'
$ws.Range("F2").Copy()
$ws.Range("F2").PasteSpecial(-4163)
$ws.Range("G2").Value = "FÍSICA MECANICA / 3"
